# Update the "Contoso Chai Tea market trends 2023" header row:
#  - Capitalize "chai" -> "Chai" in several header labels
#  - Rename a couple of header labels
#  - Make the affected header run text bold (it already renders bold via
#    the paragraph mark formatting, but the run-level w:b was explicitly
#    set to false; flip it to true to match)

$d = $word.ActiveDocument

# Map of old header text -> new header text (each is unique in the doc).
$replacements = @{
    "Total de vendas de chai (unidades)"          = "Total de vendas de Chai (unidades)";
    "Vendas de chai artesanal (unidades)"         = "Vendas de Chai artesanal (unidades)";
    "Vendas de chai pronto (unidades)"            = "Vendas de Chai pré-fabricado (unidades)";
    "Participação nas redes sociais (exibições)"  = "Engajamento nas redes sociais (visualizações)";
    "Pesquisas online por chai"                   = "Pesquisas online por Chai";
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]

    $range = $d.Content
    [void]$range.Find.Execute($old, $true, $true, $false, $false, $false, `
                               $true, 1, $false, $new, 2)
}

# Now ensure each new header label's run is bold.
foreach ($new in $replacements.Values) {
    $range = $d.Content
    $found = $range.Find.Execute($new, $true, $true, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if ($found) {
        $range.Font.Bold = 1
    }
}
